$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "83×11=913";  New = "86×49=4214" },
    @{ Old = "69×14=966";  New = "17×77=1309" },
    @{ Old = "40×51=2040"; New = "16×48=768" },
    @{ Old = "25×98=2450"; New = "60×32=1920" },
    @{ Old = "28×80=2240"; New = "46×79=3634" },
    @{ Old = "41×40=1640"; New = "81×73=5913" },
    @{ Old = "17×13=221";  New = "63×88=5544" },
    @{ Old = "65×58=3770"; New = "76×95=7220" },
    @{ Old = "88×26=2288"; New = "57×46=2622" },
    @{ Old = "92×19=1748"; New = "48×23=1104" },
    @{ Old = "67×25=1675"; New = "20×69=1380" },
    @{ Old = "38×31=1178"; New = "21×16=336" },
    @{ Old = "98×84=8232"; New = "53×13=689" },
    @{ Old = "93×97=9021"; New = "57×70=3990" },
    @{ Old = "38×29=1102"; New = "67×53=3551" },
    @{ Old = "88×42=3696"; New = "63×36=2268" },
    @{ Old = "75×55=4125"; New = "74×84=6216" },
    @{ Old = "90×19=1710"; New = "98×57=5586" },
    @{ Old = "64×69=4416"; New = "72×36=2592" },
    @{ Old = "72×18=1296"; New = "11×96=1056" },
    @{ Old = "88×47=4136"; New = "49×27=1323" },
    @{ Old = "60×31=1860"; New = "79×19=1501" },
    @{ Old = "52×30=1560"; New = "14×28=392" },
    @{ Old = "58×30=1740"; New = "25×79=1975" },
    @{ Old = "51×84=4284"; New = "99×22=2178" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
